# Weekly fruit/hortaliza update: insert a new weekly (Primera/Segunda) pair of
# rows at the top of the date-ordered block (old row 108), pushing the
# existing rows 108-167 down to 110-169, and populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 108; Excel shifts rows 108:167 down to 110:169
# and the whole used range grows from R167 to R169.
$ws.Rows("108:109").Insert()

# New row 108 ("Primera" quality) - $/atado 0,5 a 1 kilo, Región de Ñuble
$ws.Cells.Item(108, 1).Value = 11
$ws.Cells.Item(108, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(108, 3).Value = "Bíobío"
$ws.Cells.Item(108, 4).Value = 44825
$ws.Cells.Item(108, 5).Value = 8
$ws.Cells.Item(108, 6).Value = 100112044
$ws.Cells.Item(108, 7).Value = "Perejil"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 200
$ws.Cells.Item(108, 11).Value = 700
$ws.Cells.Item(108, 12).Value = 800
$ws.Cells.Item(108, 13).Value = 750
$ws.Cells.Item(108, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(108, 15).Value = "Región de Ñuble"
$ws.Cells.Item(108, 16).Value = 750
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = "Hortaliza"

# New row 109 ("Segunda" quality) - $/atado 0,5 a 1 kilo, Región de Ñuble
$ws.Cells.Item(109, 1).Value = 11
$ws.Cells.Item(109, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(109, 3).Value = "Bíobío"
$ws.Cells.Item(109, 4).Value = 44825
$ws.Cells.Item(109, 5).Value = 8
$ws.Cells.Item(109, 6).Value = 100112044
$ws.Cells.Item(109, 7).Value = "Perejil"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Segunda"
$ws.Cells.Item(109, 10).Value = 100
$ws.Cells.Item(109, 11).Value = 600
$ws.Cells.Item(109, 12).Value = 600
$ws.Cells.Item(109, 13).Value = 600
$ws.Cells.Item(109, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(109, 15).Value = "Región de Ñuble"
$ws.Cells.Item(109, 16).Value = 600
$ws.Cells.Item(109, 17).Value = 1
$ws.Cells.Item(109, 18).Value = "Hortaliza"

# Ensure the date cells keep the same date/time number format used elsewhere
# in column D.
$ws.Cells.Item(108, 4).NumberFormat = $ws.Cells.Item(110, 4).NumberFormat
$ws.Cells.Item(109, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat
